$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the real closing date for rows 7 and 8
$ws.Range("E7").Value = (Get-Date -Year 2016 -Month 3 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E7").NumberFormat = "DD/MM/YY"

$ws.Range("E8").Value = (Get-Date -Year 2016 -Month 3 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E8").NumberFormat = "DD/MM/YY"

# Update status to "Cerrada" (Closed)
$ws.Range("F7").Value = "Cerrada"
$ws.Range("F8").Value = "Cerrada"

# Update the active selection to F8
$ws.Range("F8").Select()
